$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (I1, J1) - copy style/format from existing header cell (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-8: I column is always 1, J column mirrors H column
$hValues = @{
    2 = 4
    3 = 7
    4 = 6
    5 = 5
    6 = 5
    7 = 4
    8 = 2
}

foreach ($r in 2..8) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hValues[$r]
}
